{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, large\n// numbers) in specific resume bullet paragraphs by splitting the paragraph's\n// single run into multiple runs, applying bold + color (#2C3E50) to the\n// numeric spans while leaving the surrounding text unformatted.\n//\n// The paragraphs touched are matched by their exact current text so the\n// script is resilient to paragraph re-ordering / indexing differences.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Each plan: { before: <exact current paragraph text>, segments: [...] }\n// A segment is { text, hl } where hl === true means \"bold + colored\".\nconst plans = [\n  {\n    before:\n      \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    segments: [\n      {\n        text:\n          \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from \",\n      },\n      { text: \"23%\", hl: true },\n      { text: \" to \" },\n      { text: \"64%\", hl: true },\n    ],\n  },\n  {\n    before:\n      \"\\u2022 Utilized advanced sampling methods to decrease survey margin of error from \\u00B14.2% to \\u00B12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\",\n    segments: [\n      {\n        text:\n          \"\\u2022 Utilized advanced sampling methods to decrease survey margin of error from \",\n      },\n      { text: \"\\u00B14.2%\", hl: true },\n      { text: \" to \" },\n      { text: \"\\u00B12.1%\", hl: true },\n      { text: \", increasing voter turnout prediction accuracy from \" },\n      { text: \"71%\", hl: true },\n      { text: \" to \" },\n      { text: \"87%\", hl: true },\n      {\n        text:\n          \", and ensuring survey results more closely reflected true population attitudes\",\n      },\n    ],\n  },\n  {\n    before:\n      \"\\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    segments: [\n      {\n        text:\n          \"\\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by \",\n      },\n      { text: \"73.5%\", hl: true },\n      { text: \", saving campaigns and organizations \" },\n      { text: \"$4.7M\", hl: true },\n      { text: \" and enabling smaller nonprofits to conduct analysis\" },\n    ],\n  },\n  {\n    before:\n      \"\\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n    segments: [\n      {\n        text:\n          \"\\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over \",\n      },\n      { text: \"$2\", hl: true },\n      { text: \" trillion\" },\n    ],\n  },\n  {\n    before:\n      \"\\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%\",\n    segments: [\n      {\n        text:\n          \"\\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by \",\n      },\n      { text: \"57%\", hl: true },\n    ],\n  },\n  {\n    before:\n      \"\\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \\u00B14.2% to \\u00B12.1%\",\n    segments: [\n      {\n        text:\n          \"\\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \",\n      },\n      { text: \"\\u00B14.2%\", hl: true },\n      { text: \" to \" },\n      { text: \"\\u00B12.1%\", hl: true },\n    ],\n  },\n  {\n    before: \"\\u2022 Increased voter turnout prediction accuracy from 71% to 87%\",\n    segments: [\n      { text: \"\\u2022 Increased voter turnout prediction accuracy from \" },\n      { text: \"71%\", hl: true },\n      { text: \" to \" },\n      { text: \"87%\", hl: true },\n    ],\n  },\n  {\n    before:\n      \"\\u2022 Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\",\n    segments: [\n      {\n        text:\n          \"\\u2022 Methodological advancement: Improved segmentation accuracy \",\n      },\n      { text: \"34%\", hl: true },\n      { text: \" and survey incidence \" },\n      { text: \"28%\", hl: true },\n    ],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Track which plans have already been consumed so that, in the (unlikely)\n// case of duplicate paragraph text, each plan is only applied once per match\n// found, left to right.\nconst used = new Array(plans.length).fill(false);\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n\n  for (let p = 0; p < plans.length; p++) {\n    if (text === plans[p].before) {\n      // Rebuild this paragraph's runs from scratch.\n      para.clear();\n      await context.sync();\n\n      for (const seg of plans[p].segments) {\n        const r = para.insertText(seg.text, \"End\");\n        if (seg.hl) {\n          r.font.bold = true;\n          r.font.color = HIGHLIGHT_COLOR;\n        }\n        await context.sync();\n      }\n      break;\n    }\n  }\n}\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, large\n# numbers) in specific resume bullet paragraphs. Each target paragraph's\n# single run of text is effectively split into multiple runs by applying\n# Bold + Color (#2C3E50) formatting to Range objects covering just the\n# numeric substrings; Word keeps the surrounding text in separate,\n# unformatted runs.\n\n$doc = $word.ActiveDocument\n\n$PM = [char]0xB1   # \"\u00b1\"\n\nfunction ColorVal([int]$r, [int]$g, [int]$b) {\n    # Word COM Font.Color wants a single int packed as 0x00BBGGRR.\n    return $r + ($g * 256) + ($b * 65536)\n}\n$HighlightColor = ColorVal 0x2C 0x3E 0x50\n\n# Build the bullet character once and reuse it (avoids unicode escape\n# ambiguity in the literal strings below).\n$Bullet = [char]0x2022\n\n$plans = @(\n    @{\n        Text    = \"$Bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Targets = @(\"23%\", \"64%\")\n    },\n    @{\n        Text    = \"$Bullet Utilized advanced sampling methods to decrease survey margin of error from ${PM}4.2% to ${PM}2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\"\n        Targets = @(\"${PM}4.2%\", \"${PM}2.1%\", \"71%\", \"87%\")\n    },\n    @{\n        Text    = \"$Bullet Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\"\n        Targets = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        Text    = \"$Bullet Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\"\n        Targets = @(\"`$2\")\n    },\n    @{\n        Text    = \"$Bullet Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%\"\n        Targets = @(\"57%\")\n    },\n    @{\n        Text    = \"$Bullet Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ${PM}4.2% to ${PM}2.1%\"\n        Targets = @(\"${PM}4.2%\", \"${PM}2.1%\")\n    },\n    @{\n        Text    = \"$Bullet Increased voter turnout prediction accuracy from 71% to 87%\"\n        Targets = @(\"71%\", \"87%\")\n    },\n    @{\n        Text    = \"$Bullet Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\"\n        Targets = @(\"34%\", \"28%\")\n    }\n)\n\n$count = $doc.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $doc.Paragraphs.Item($i)\n    $pr = $para.Range\n    $full = $pr.Text.TrimEnd([char]13, [char]7)\n\n    foreach ($plan in $plans) {\n        if ($full -eq $plan.Text) {\n            $pStart = $pr.Start\n            $searchFrom = 0\n            foreach ($needle in $plan.Targets) {\n                $idx = $full.IndexOf($needle, $searchFrom)\n                if ($idx -ge 0) {\n                    $absStart = $pStart + $idx\n                    $absEnd = $absStart + $needle.Length\n                    $sub = $doc.Range($absStart, $absEnd)\n                    $sub.Font.Bold = 1\n                    $sub.Font.Color = $HighlightColor\n                    $searchFrom = $idx + $needle.Length\n                }\n            }\n            break\n        }\n    }\n}\n"}
